$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 37

# Date / Week-number columns look numeric to Excel's smart-entry parser
# (ISO date, zero-padded week number) - force them to Text first so the
# literal strings are preserved exactly, then drop the Text number format
# again so the cell is left with the sheet's normal (default) styling.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value  = "2024-01-09"
$ws.Cells.Item($row, 1).ClearFormats()
$ws.Cells.Item($row, 2).Value  = "09:12:37"
$ws.Cells.Item($row, 3).Value  = "Tuesday"
$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value  = "01"
$ws.Cells.Item($row, 4).ClearFormats()
$ws.Cells.Item($row, 5).Value  = 139490
$ws.Cells.Item($row, 6).Value  = 142574
$ws.Cells.Item($row, 7).Value  = 171869
$ws.Cells.Item($row, 8).Value  = 147227
$ws.Cells.Item($row, 9).Value  = -1
$ws.Cells.Item($row, 10).Value = 117870
$ws.Cells.Item($row, 11).Value = 224600
$ws.Cells.Item($row, 12).Value = 249808
$ws.Cells.Item($row, 13).Value = 185086
$ws.Cells.Item($row, 14).Value = 110383
$ws.Cells.Item($row, 15).Value = 40631
$ws.Cells.Item($row, 16).Value = 30812
$ws.Cells.Item($row, 17).Value = 72388
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 41460
$ws.Cells.Item($row, 20).Value = -1
